$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.526.36"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "3.059.02"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'385.94"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").Value = "'102.98"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").Value = "'0.543"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.585"
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("D10").Value = "'36.84"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "'0.0860"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "3.545.13"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").Value = "'18.55"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "3.054.38"
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("D17").Value = "'0.973"
$ws.Range("E17").Value = "  -2.65%  "
$ws.Range("D18").Value = "'10.65"
$ws.Range("E18").Value = "  -4.33%  "
$ws.Range("D19").Value = "51.575.08"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "'3.15"
$ws.Range("E20").Value = "  +2.04%  "
$ws.Range("D21").Value = "'12.44"
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("D22").Value = "0.0₃0967"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").Value = "'70.17"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").Value = "'268.08"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "'3.14"
$ws.Range("E25").Value = "  -2.35%  "
$ws.Range("D26").Value = "'8.21"
$ws.Range("E26").Value = "  +4.18%  "
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'7.26"
$ws.Range("E28").Value = "  -2.58%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.170"
$ws.Range("E29").Value = "  +1.57%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -1.70%  "
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("D33").Value = "'34.72"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  +0.98%  "
$ws.Range("E35").Value = "  -3.03%  "
$ws.Range("D36").Value = "'0.0447"
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "'3.32"
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("D39").Value = "'0.293"
$ws.Range("E39").Value = "  +7.78%  "
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("E41").Value = "  +1.36%  "
$ws.Range("D42").Value = "'2.57"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").Value = "'125.28"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").Value = "'3.75"
$ws.Range("E45").Value = "  +2.80%  "
$ws.Range("D46").Value = "'21.89"
$ws.Range("E46").Value = "  +1.73%  "
$ws.Range("D47").Value = "'2.09"
$ws.Range("E47").Value = "  +2.98%  "
$ws.Range("D48").Value = "'2.44"
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("D49").Value = "2.034.09"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").Value = "3.359.11"
$ws.Range("E50").Value = "  +2.00%  "
$ws.Range("E51").Value = "  +6.49%  "
